$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Three new trailing columns of header data
$ws.Range("D1").Value = "ORG_FAC_IDENOLD"
$ws.Range("E1").Value = "ORG_FAC_IDENNEW"
$ws.Range("F1").Value = "ORG_FAC_STATUS"

# Re-balance column widths: A:C / G:K a touch wider, the new D:F narrower
$ws.Range("A1:C1").EntireColumn.ColumnWidth = 34.6
$ws.Range("D1:F1").EntireColumn.ColumnWidth = 23.6
$ws.Range("G1:K1").EntireColumn.ColumnWidth = 34.6

# Move the active selection to F6, matching the saved view
$ws.Range("F6").Select() | Out-Null
